$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '28.051.17'
$ws.Range("E2").Value = '  -0.11%  '

$ws.Range("D3").Value = '1.809.08'
$ws.Range("E3").Value = '  +1.12%  '

Set-TextValue "D4" '1.007'
$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue "D5" '1.009'
$ws.Range("E5").Value = '  +0.65%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue "D6" '307.70'
$ws.Range("E6").Value = '  -1.86%  '

Set-TextValue "D7" '0.5003'
$ws.Range("E7").Value = '  -3.72%  '

Set-TextValue "D8" '0.3868'
$ws.Range("E8").Value = '  +1.87%  '

Set-TextValue "D9" '0.09461'
$ws.Range("E9").Value = '  +18.82%  '

Set-TextValue "D10" '1.099'
$ws.Range("E10").Value = '  +0.97%  '

$ws.Range("E11").Value = '  -1.49%  '

Set-TextValue "D12" '6.377'
$ws.Range("E12").Value = '  +2.03%  '

Set-TextValue "D13" '1.007'
$ws.Range("E13").Value = '  +0.35%  '

Set-TextValue "D14" '20.60'
$ws.Range("E14").Value = '  +0.69%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.834.30'
$ws.Range("E15").Value = '  +2.43%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D16" '7.246'
$ws.Range("E16").Value = '  -0.33%  '

Set-TextValue "D17" '0.00001120'
$ws.Range("E17").Value = '  +2.83%  '

Set-TextValue "D18" '92.48'
$ws.Range("E18").Value = '  +1.47%  '

Set-TextValue "D19" '0.06598'
$ws.Range("E19").Value = '  +0.86%  '

Set-TextValue "D20" '1.009'
$ws.Range("E20").Value = '  +0.55%  '

Set-TextValue "D21" '17.12'
$ws.Range("E21").Value = '  -0.81%  '

Set-TextValue "D22" '5.972'
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").Value = '28.085.76'
$ws.Range("E23").Value = '  -0.18%  '

Set-TextValue "D24" '11.07'
$ws.Range("E24").Value = '  -0.15%  '

Set-TextValue "D25" '2.237'
$ws.Range("E25").Value = '  -0.98%  '

Set-TextValue "D26" '157.76'
$ws.Range("E26").Value = '  -1.22%  '

$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '2.031.62'
$ws.Range("E27").Value = '  +1.92%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D28" '20.56'
$ws.Range("E28").Value = '  +1.08%  '

Set-TextValue "D29" '2.380'
$ws.Range("E29").Value = '  +2.45%  '

Set-TextValue "D30" '127.43'
$ws.Range("E30").Value = '  +4.14%  '

Set-TextValue "D31" '0.1081'
$ws.Range("E31").Value = '  -0.14%  '

Set-TextValue "D32" '1.047'
$ws.Range("E32").Value = '  -0.47%  '

$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D33" '3.644'
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D34" '5.547'
$ws.Range("E34").Value = '  +0.49%  '

Set-TextValue "D35" '0.06837'
$ws.Range("E35").Value = '  -4.73%  '

Set-TextValue "D36" '8.937'
$ws.Range("E36").Value = '  +4.44%  '

Set-TextValue "D37" '0.02313'
$ws.Range("E37").Value = '  +0.43%  '

Set-TextValue "D38" '0.2138'
$ws.Range("E38").Value = '  +0.22%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D39" '4.964'
$ws.Range("E39").Value = '  -1.95%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D40" '11.37'
$ws.Range("E40").Value = '  -6.30%  '

Set-TextValue "D41" '0.6155'
$ws.Range("E41").Value = '  +0.17%  '

Set-TextValue "D42" '1.007'
$ws.Range("E42").Value = '  +0.49%  '

Set-TextValue "D43" '1.151'
$ws.Range("E43").Value = '  -0.62%  '

Set-TextValue "D44" '13.00'
$ws.Range("E44").Value = '  -1.69%  '

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue "D45" '1.296'
$ws.Range("E45").Value = '  -5.35%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D46" '0.5905'
$ws.Range("E46").Value = '  +0.01%  '

Set-TextValue "D47" '3.691'
$ws.Range("E47").Value = '  -1.73%  '

Set-TextValue "D48" '124.18'
$ws.Range("E48").Value = '  -1.92%  '

Set-TextValue "D49" '1.957'
$ws.Range("E49").Value = '  +2.32%  '

Set-TextValue "D50" '1.171'
$ws.Range("E50").Value = '  -3.22%  '

Set-TextValue "D51" '0.06748'
$ws.Range("E51").Value = '  -0.14%  '
